$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.038.84"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "2.382.10"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.71"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.92"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.64"
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "2.746.38"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "2.380.89"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.827"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "45.937.63"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -5.40%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.06"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.96"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.88"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -5.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "39.83"
$ws.Range("E27").Value = "  -11.08%  "
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.82"
$ws.Range("E30").Value = "  +21.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.95"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.80"
$ws.Range("E32").Value = "  +6.88%  "
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "146.88"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0776"
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +6.64%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.13"
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.24"
$ws.Range("E42").Value = "  -8.18%  "
$ws.Range("D43").Value = "1.937.56"
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.79"
$ws.Range("E45").Value = "  +4.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.79"
$ws.Range("E46").Value = "  -9.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.47"
$ws.Range("E47").Value = "  +4.83%  "
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.00"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "2.615.84"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.07"
$ws.Range("E51").Value = "  -7.46%  "
